# "Updated with data from Apr 9"
#
# The source CSV this sheet tracks gained a new daily row, and the
# previous day's "tested" figure (which hadn't been published yet) became
# available and was back-filled.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Back-fill the "tested" count for row 34 (date 2020-04-07), which was
# left blank when that row was first entered.
$ws.Range("B34").Value = 842

# Append row 35 (date 2020-04-08) with the day's figures. Columns not
# listed here (E, G, I, K, M, O, Q, S, U, W, BF, ...) have no data for
# this date and are left blank, same as in the preceding rows.
$row35 = @{
    "A"  = 43929
    "C"  = 1530
    "D"  = 10
    "F"  = 16
    "H"  = 244
    "J"  = 301
    "L"  = 268
    "N"  = 278
    "P"  = 201
    "R"  = 123
    "T"  = 86
    "V"  = 3
    "X"  = 737
    "Y"  = 785
    "Z"  = 8
    "AA" = 316
    "AB" = 122
    "AC" = 36
    "AD" = 43
    "AE" = 133
    "AF" = 4
    "AG" = 8
    "AH" = 88
    "AI" = 32
    "AJ" = 35
    "AK" = 9
    "AL" = 29
    "AM" = 14
    "AN" = 28
    "AO" = 34
    "AP" = 16
    "AQ" = 783
    "AR" = 20
    "AS" = 16
    "AT" = 5
    "AU" = 24
    "AV" = 1
    "AW" = 13
    "AX" = 1
    "AY" = 1
    "AZ" = 6
    "BA" = 3
    "BB" = 14
    "BC" = 2
    "BD" = 10
    "BE" = 14
    "BG" = 37
    "BH" = 3
    "BI" = 85
}

foreach ($col in $row35.Keys) {
    $ws.Range($col + "35").Value = $row35[$col]
}

# Leave the workbook's selection/scroll state the way it was left after
# entering the new row (cursor parked one cell past the last data column).
$ws.Range("BJ35").Select()
$excel.ActiveWindow.ScrollColumn = 50
$excel.ActiveWindow.ScrollRow = 1
